$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Fill-Cell($cell, $paraId, $textId, $text) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="' + $paraId + '" w14:textId="' + $textId + '" w:rsidR="00685ECC" w:rsidRPr="00A6112F" w:rsidRDefault="00685ECC" w:rsidP="00685ECC">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:outlineLvl w:val="0"/>' +
        '<w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="2D3B45"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="2D3B45"/></w:rPr><w:t>' + $text + '</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xml)
}

Fill-Cell $t.Cell(5, 1) "142B03E6" "61F2E110" "19"
Fill-Cell $t.Cell(5, 2) "51413AD6" "3DC5A25F" "Hunter Malinowski"
Fill-Cell $t.Cell(5, 3) "4E7B0AFE" "0AECE97A" "Developer"
